$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets("CAN1")

# New column E header first (so "Attributes.ReadRate" is appended to the
# shared-string table before "Attributes.Baudrate"), then the D column
# header text is changed from "Attributes" to "Attributes.Baudrate".
$ws.Range("E1").Value = "Attributes.ReadRate"
$ws.Range("D1").Value = "Attributes.Baudrate"

# Row 2 values become real numbers (variant) instead of text.
$ws.Range("D2").Value = 500000
$ws.Range("E2").Value = 10.5

# Give the new column a width similar to the other data columns.
$ws.Columns.Item(5).ColumnWidth = 28.5

# Make CAN1 the active sheet / tab, with E3 selected.
$ws.Activate()
$ws.Range("E3").Select()
